# Update the date header and each arithmetic-expression cell in the table
# to the new values from the commit (old -> new text pairs).
$d = $word.ActiveDocument

$pairs = @(
    @("2025-09-16 Tuesday", "2025-09-17 Wednesday"),
    @("99-53=", "85+12="),
    @("83-35=", "36+10="),
    @("71+28=", "45-20="),
    @("43+46=", "89+2="),
    @("49+38=", "63-0="),
    @("10+35=", "55+7="),
    @("34+4=", "99-41="),
    @("5+94=", "80-3="),
    @("88-71=", "10+81="),
    @("70-10=", "94-91="),
    @("12+20=", "45-11="),
    @("61+36=", "17+28="),
    @("81-14=", "42-9="),
    @("2+19=", "34+50="),
    @("67-17=", "33+39="),
    @("48+34=", "82-19="),
    @("82-69=", "81-68="),
    @("37-1=", "51+10="),
    @("76-70=", "16+82="),
    @("20+41=", "71-67="),
    @("42+15=", "2+26="),
    @("22+74=", "72-51="),
    @("45+0=", "93+5="),
    @("68-55=", "13+27="),
    @("46-43=", "96-46="),
    @("78+16=", "15-15="),
    @("1+95=", "36-19="),
    @("72+25=", "41-35="),
    @("87-31=", "70-34="),
    @("68+18=", "33-0="),
    @("6+18=", "30+36="),
    @("86-45=", "56+10="),
    @("39+21=", "5+44="),
    @("5-2=", "49+44="),
    @("38+59=", "57-13="),
    @("41-5=", "66-17="),
    @("65-58=", "9+30="),
    @("48-47=", "29+60="),
    @("78-58=", "63+2="),
    @("24+8=", "2+84="),
    @("94-55=", "19+22="),
    @("72-60=", "37-35="),
    @("70-61=", "80-74="),
    @("15+20=", "90-21="),
    @("34+7=", "51-22="),
    @("13+19=", "76-2="),
    @("73+3=", "97-23="),
    @("98-65=", "79-63="),
    @("83-46=", "22+47="),
    @("68-44=", "98-81="),
    @("0+8=", "28+18="),
    @("88-76=", "56-18="),
    @("97-2=", "74-48="),
    @("16+65=", "52-13="),
    @("0+47=", "32+8="),
    @("65+17=", "92-0="),
    @("1+96=", "21+54="),
    @("93-79=", "51-41="),
    @("7+22=", "13+60="),
    @("95-55=", "59+29="),
    @("6+7=", "49-45="),
    @("47-20=", "96-95="),
    @("27+1=", "65-21="),
    @("80-70=", "38+25="),
    @("11+50=", "65+23="),
    @("92+6=", "73+5="),
    @("46-5=", "55-0="),
    @("77-59=", "64-18="),
    @("43+49=", "73+21="),
    @("81-23=", "41-1="),
    @("42+34=", "40+12="),
    @("88-21=", "51+44="),
    @("84-16=", "55-54="),
    @("29+57=", "26-16="),
    @("22+65=", "41+47="),
    @("85-11=", "95-1="),
    @("22+67=", "30-11="),
    @("86+13=", "14+77="),
    @("13+68=", "3+8="),
    @("59+10=", "99-82="),
    @("72+21=", "25+16="),
    @("65+9=", "57-10="),
    @("70-44=", "61+6="),
    @("44-43=", "72-13="),
    @("65+7=", "98-35="),
    @("53-51=", "52+7="),
    @("23-11=", "98-92="),
    @("34+16=", "27+19="),
    @("22-11=", "67-4="),
    @("98-91=", "26+19="),
    @("49+28=", "28+55="),
    @("32-14=", "52+46="),
    @("48+25=", "61+15="),
    @("8+16=", "11+26="),
    @("30+17=", "6+53="),
    @("90-33=", "50-33="),
    @("63+34=", "43-10="),
    @("88-44=", "89-25="),
    @("16+30=", "33-26="),
    @("53-16=", "96-83=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}
